$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 40000
$ws.Range("J3").Value = 40000
$ws.Range("L3").Value = 40000
$ws.Range("N3").Value = -40228
$ws.Range("H39").Value = 495.6875
$ws.Range("J39").Value = 1159
$ws.Range("L39").Value = 3477
$ws.Range("N39").Value = -4069
$ws.Range("H100").Value = 9667.666999999999
$ws.Range("I100").Value = 4998.5
$ws.Range("K100").Value = 4998.5
$ws.Range("M100").Value = -4457.5
$ws.Range("H102").Value = 40000
$ws.Range("J102").Value = 40000
$ws.Range("L102").Value = 40000
$ws.Range("N102").Value = -46490
$ws.Range("H137").Value = 3077.5293
$ws.Range("I137").Value = 2025.5
$ws.Range("J137").Value = 3456.26
$ws.Range("K137").Value = 6076.5
$ws.Range("L137").Value = 10368.78
$ws.Range("M137").Value = -3526.5
$ws.Range("N137").Value = -15468.78
$ws.Range("H138").Value = 2606.5557
$ws.Range("I138").Value = 1238.8928
$ws.Range("J138").Value = 3476.8865
$ws.Range("K138").Value = 3716.6784
$ws.Range("L138").Value = 10430.6595
$ws.Range("M138").Value = 1423.3216
$ws.Range("N138").Value = -20710.6595

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 33959.668
$ws.Range("I2").Value = 939.5
$ws.Range("K2").Value = 939.5
$ws.Range("M2").Value = -826.5
$ws.Range("H28").Value = 5466
$ws.Range("I28").Value = 5466
$ws.Range("K28").Value = 5466
$ws.Range("M28").Value = -5274
$ws.Range("H45").Value = 1566.3103
$ws.Range("I45").Value = 1573.8695
$ws.Range("K45").Value = 1573.8695
$ws.Range("M45").Value = -1196.8695
$ws.Range("H74").Value = 10419547
$ws.Range("I74").Value = 11906809
$ws.Range("J74").Value = 8710.5
$ws.Range("K74").Value = 11906809
$ws.Range("L74").Value = 8710.5
$ws.Range("M74").Value = -11905935
$ws.Range("N74").Value = -10458.5
$ws.Range("H77").Value = 10419547
$ws.Range("I77").Value = 11906809
$ws.Range("J77").Value = 8710.5
$ws.Range("K77").Value = 59534045
$ws.Range("L77").Value = 43552.5
$ws.Range("M77").Value = -59529677
$ws.Range("N77").Value = -52288.5
$ws.Range("H99").Value = 5466
$ws.Range("I99").Value = 5466
$ws.Range("K99").Value = 5466
$ws.Range("M99").Value = -2471
$ws.Range("H102").Value = 1601.6666
$ws.Range("I102").Value = 1586.25
$ws.Range("J102").Value = 1725
$ws.Range("K102").Value = 1586.25
$ws.Range("L102").Value = 1725
$ws.Range("M102").Value = 35.75
$ws.Range("N102").Value = -4969
$ws.Range("H109").Value = 91732.336
$ws.Range("J109").Value = 91732.336
$ws.Range("L109").Value = 91732.336
$ws.Range("N109").Value = -94506.336
$ws.Range("H116").Value = 33959.668
$ws.Range("I116").Value = 939.5
$ws.Range("K116").Value = 939.5
$ws.Range("M116").Value = 1354.5
$ws.Range("H132").Value = 4136.0356
$ws.Range("I132").Value = 2133.0476
$ws.Range("J132").Value = 10145
$ws.Range("K132").Value = 6399.1428
$ws.Range("L132").Value = 30435
$ws.Range("M132").Value = -3869.1428
$ws.Range("N132").Value = -35495

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 33959.668
$ws.Range("I3").Value = 939.5
$ws.Range("K3").Value = 939.5
$ws.Range("M3").Value = -825.5
$ws.Range("H80").Value = 532.6111
$ws.Range("I80").Value = 441.8
$ws.Range("J80").Value = 567.53845
$ws.Range("K80").Value = 441.8
$ws.Range("L80").Value = 567.53845
$ws.Range("M80").Value = 556.2
$ws.Range("N80").Value = -2563.53845
$ws.Range("H83").Value = 532.6111
$ws.Range("I83").Value = 441.8
$ws.Range("J83").Value = 567.53845
$ws.Range("K83").Value = 2209
$ws.Range("L83").Value = 2837.69225
$ws.Range("M83").Value = 2783
$ws.Range("N83").Value = -12821.69225
$ws.Range("H94").Value = 1812.56
$ws.Range("I94").Value = 1168.3334
$ws.Range("K94").Value = 1168.3334
$ws.Range("M94").Value = -717.3334
$ws.Range("H105").Value = 29449
$ws.Range("I105").Value = 100995.5
$ws.Range("J105").Value = 11562.375
$ws.Range("K105").Value = 100995.5
$ws.Range("L105").Value = 11562.375
$ws.Range("M105").Value = -99248.5
$ws.Range("N105").Value = -15056.375
$ws.Range("H107").Value = 1321.7941
$ws.Range("I107").Value = 1082.8928
$ws.Range("J107").Value = 2436.6667
$ws.Range("K107").Value = 1082.8928
$ws.Range("L107").Value = 2436.6667
$ws.Range("M107").Value = 837.1071999999999
$ws.Range("N107").Value = -6276.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1636.6923
$ws.Range("I16").Value = 1152.4546
$ws.Range("K16").Value = 1152.4546
$ws.Range("M16").Value = -865.4546
$ws.Range("H86").Value = 8749
$ws.Range("I86").Value = 7000
$ws.Range("K86").Value = 7000
$ws.Range("M86").Value = -5877
$ws.Range("H89").Value = 8749
$ws.Range("I89").Value = 7000
$ws.Range("K89").Value = 35000
$ws.Range("M89").Value = -29384
$ws.Range("H105").Value = 3774.875
$ws.Range("I105").Value = 1061
$ws.Range("K105").Value = 1061
$ws.Range("M105").Value = 686
$ws.Range("H113").Value = 1636.6923
$ws.Range("I113").Value = 1152.4546
$ws.Range("K113").Value = 1152.4546
$ws.Range("M113").Value = 1017.5454
$ws.Range("H132").Value = 2686.3103
$ws.Range("I132").Value = 2014.4231
$ws.Range("K132").Value = 6043.2693
$ws.Range("M132").Value = -3513.2693

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 5915997.5
$ws.Range("I4").Value = 4166987
$ws.Range("J4").Value = 9144940
$ws.Range("K4").Value = 12500961
$ws.Range("L4").Value = 27434820
$ws.Range("M4").Value = -12500849
$ws.Range("N4").Value = -27435044
$ws.Range("H107").Value = 1954051
$ws.Range("I107").Value = 641.75
$ws.Range("J107").Value = 2605187.2
$ws.Range("K107").Value = 1925.25
$ws.Range("L107").Value = 7815561.600000001
$ws.Range("M107").Value = -5.25
$ws.Range("N107").Value = -7819401.600000001
$ws.Range("H121").Value = 775
$ws.Range("I121").Value = 660
$ws.Range("K121").Value = 1980
$ws.Range("M121").Value = -670
$ws.Range("H131").Value = 6521092.5
$ws.Range("J131").Value = 4987609.5
$ws.Range("L131").Value = 14962828.5
$ws.Range("N131").Value = -14972908.5
$ws.Range("H132").Value = 4409.273
$ws.Range("I132").Value = 3063.125
$ws.Range("J132").Value = 7999
$ws.Range("K132").Value = 27568.125
$ws.Range("L132").Value = 71991
$ws.Range("M132").Value = -25038.125
$ws.Range("N132").Value = -77051

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 13293.5
$ws.Range("I70").Value = 9771.154
$ws.Range("J70").Value = 17456.273
$ws.Range("K70").Value = 9771.154
$ws.Range("L70").Value = 17456.273
$ws.Range("M70").Value = -9501.154
$ws.Range("N70").Value = -17996.273
$ws.Range("H73").Value = 13293.5
$ws.Range("I73").Value = 9771.154
$ws.Range("J73").Value = 17456.273
$ws.Range("K73").Value = 9771.154
$ws.Range("L73").Value = 17456.273
$ws.Range("M73").Value = -8835.154
$ws.Range("N73").Value = -19328.273
$ws.Range("H80").Value = 337966.4
$ws.Range("J80").Value = 6624.25
$ws.Range("L80").Value = 6624.25
$ws.Range("N80").Value = -8620.25
$ws.Range("H83").Value = 337966.4
$ws.Range("J83").Value = 6624.25
$ws.Range("L83").Value = 33121.25
$ws.Range("N83").Value = -43105.25
$ws.Range("H122").Value = 4135.8125
$ws.Range("I122").Value = 3415.258
$ws.Range("J122").Value = 5449.7646
$ws.Range("K122").Value = 10245.774
$ws.Range("L122").Value = 16349.2938
$ws.Range("M122").Value = -7795.773999999999
$ws.Range("N122").Value = -21249.2938
$ws.Range("H126").Value = 3615.0344
$ws.Range("I126").Value = 2489.2666
$ws.Range("K126").Value = 7467.7998
$ws.Range("M126").Value = -4997.7998

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6232.6665
$ws.Range("I7").Value = 5215.7383
$ws.Range("K7").Value = 5215.7383
$ws.Range("M7").Value = -5103.7383
$ws.Range("H55").Value = 2274705.2
$ws.Range("I55").Value = 3333907.8
$ws.Range("K55").Value = 3333907.8
$ws.Range("M55").Value = -3333734.8
$ws.Range("H93").Value = 1631.3846
$ws.Range("I93").Value = 1564.4546
$ws.Range("J93").Value = 1999.5
$ws.Range("K93").Value = 1564.4546
$ws.Range("L93").Value = 1999.5
$ws.Range("M93").Value = -316.4546
$ws.Range("N93").Value = -4495.5
$ws.Range("H100").Value = 7763.8438
$ws.Range("I100").Value = 2163
$ws.Range("J100").Value = 13364.6875
$ws.Range("K100").Value = 2163
$ws.Range("L100").Value = 13364.6875
$ws.Range("M100").Value = -1622
$ws.Range("N100").Value = -14446.6875
$ws.Range("H126").Value = 6232.6665
$ws.Range("I126").Value = 5215.7383
$ws.Range("K126").Value = 15647.2149
$ws.Range("M126").Value = -13177.2149
$ws.Range("H136").Value = 7129.087
$ws.Range("I136").Value = 2890.3635
$ws.Range("K136").Value = 8671.0905
$ws.Range("M136").Value = -6121.0905

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H56").Value = 4500
$ws.Range("J56").Value = 4500
$ws.Range("L56").Value = 4500
$ws.Range("N56").Value = -5928
$ws.Range("H100").Value = 1166.35
$ws.Range("J100").Value = 1290.7778
$ws.Range("L100").Value = 2581.5556
$ws.Range("N100").Value = -3663.5556
$ws.Range("H126").Value = 2031.8529
$ws.Range("I126").Value = 1884.963
$ws.Range("K126").Value = 5654.889
$ws.Range("M126").Value = -3184.889
